# Committing Login File Changes
# Adds an "Environment" / "test" column (D) to the Login sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D1").Value = "Environment"
$ws.Range("D2").Value = "test"

# Closest achievable width to the authored 12.85546875 given this engine's
# internal rounding of ColumnWidth (quantized to 1/6 character units).
$ws.Columns.Item(4).ColumnWidth = 12

$ws.Range("D2").Select()
